# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert a new data row (period 2509) right below the last existing
#    data row (row 19), pushing the footer rows down.
$ws.Rows.Item(20).Insert()

# Copy formatting from the row above (row 19) so the new row matches the
# existing table look (borders/fills/fonts/number formats).
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's data.
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143353580"
$ws.Range("D20").Value = "IVAN ALEXANDER CASTILLA PEÑA"
$ws.Range("E20").Value = "2509"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

# 2) Update the total "VALOR MORA" amount to reflect the new period.
$ws.Range("E11").Value = 284700

# 3) Update the "Cant. Periodos" count to reflect the new period.
$ws.Range("F13").Value = 5
